$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.652.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.21%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.866.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.42%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.31%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'326.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.17%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.4630"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.62%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3910"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.23%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07921"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.68%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9695"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.38%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'22.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.45%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.936.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +4.08%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'5.728"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.06%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'6.931"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.06952"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.12%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +1.43%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.004"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.20%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.00001006"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.07%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'16.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.23%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +0.28%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'28.654.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.17%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.320"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.43%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'11.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.73%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -1.46%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.093.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.63%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'153.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.09%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'19.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.71%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'5.716"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.87%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.56%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'119.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.94%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.09369"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.64%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.9310"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.22%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'5.329"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.76%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.344"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.16%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'3.356"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.54%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.05834"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -3.00%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.02125"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.39%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.150"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.04%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'7.905"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'0.5654"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.40%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1782"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.04%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'9.922"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.07%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.07239"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.92%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'11.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.50%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.5316"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.40%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.162"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -5.23%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.137"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -6.28%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.845"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.18%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'113.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'1.004"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.29%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.343"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.90%  "
$ws.Range("E51").Style = "Normal"
